$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F3").Value = 7
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -1
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = -2
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = 5
$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 4
$ws.Range("F14").Value = -4
$ws.Range("F15").Value = 5
$ws.Range("F16").Value = 1
$ws.Range("F17").Value = 4
$ws.Range("F19").Value = -1
$ws.Range("F20").Value = 1
$ws.Range("F21").Value = -3
$ws.Range("F22").Value = -3
$ws.Range("F23").Value = -4
$ws.Range("F24").Value = -3
$ws.Range("F26").Value = -1
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = -3
$ws.Range("F30").Value = 3
$ws.Range("F32").Value = 3
$ws.Range("F33").Value = 4
$ws.Range("F35").Value = -1
$ws.Range("F36").Value = -1
$ws.Range("F37").Value = 3
$ws.Range("F39").Value = -1
